$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233. This shifts the existing rows 233-292
# down to 234-293 (carrying all of their original data with them), and
# leaves a blank row 233 ready to be populated below.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with its data.
$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value = 44551
$ws.Cells.Item(233, 5).Value = 5
$ws.Cells.Item(233, 6).Value = 100112017
$ws.Cells.Item(233, 7).Value = "Apio"
$ws.Cells.Item(233, 8).Value = "Americana (o)"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 180
$ws.Cells.Item(233, 11).Value = 9000
$ws.Cells.Item(233, 12).Value = 9000
$ws.Cells.Item(233, 13).Value = 9000
$ws.Cells.Item(233, 14).Value = "`$/docena de matas"
$ws.Cells.Item(233, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(233, 16).Value = 1500
$ws.Cells.Item(233, 17).Value = 6
$ws.Cells.Item(233, 18).Value = "Hortaliza"

# Match the date-number style used by the other cells in column D.
$ws.Cells.Item(233, 4).NumberFormat = $ws.Cells.Item(234, 4).NumberFormat
